# Soar Tutorial Part 7 - RL: fix the "Agents" directory reference and move
# the _GoBack bookmark to sit after it (matches upstream commit that moved
# the agents-files note and renamed the directory from "Demos" to "Agents").

$d = $word.ActiveDocument

# 1) "Demo" + "s" (two italic runs) -> single italic run "Agents"
$rng = $d.Content
$rng.Find.Execute("Demos", $false, $false, $false, $false, $false, $true, 1, $false, "Agents", 2) | Out-Null

# 2) Move the _GoBack bookmark from right after "...exploration policy
#    (more on this later.)" to right after the newly-renamed "Agents" run,
#    i.e. immediately before " directory." at the end of that paragraph.
$full = $d.Range(0, $d.Content.End)
$txt = $full.Text
$idx = $txt.IndexOf("Agents directory.")
$pos = $idx + 6
$bmRange = $d.Range($pos, $pos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
